$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = -11.7492
$ws.Range("C3").Value = -11.99
$ws.Range("C5").Value = -12.3161
$ws.Range("A9").Value = -20.85699999999997
$ws.Range("C11").Value = -14.01840000000001
$ws.Range("C12").Value = -14.10050000000001
$ws.Range("A13").Value = -22.00010000000003
$ws.Range("A16").Value = -20.13819999999999
$ws.Range("A18").Value = -22.01440000000001
$ws.Range("A20").Value = -22.01270000000001
$ws.Range("C21").Value = -13.17760000000001
